$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date number format, bold, border) from the last
# existing date cell (A19) down onto the new date cell (A20), then set
# the new row's values.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -2.451276118722334
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 1.795477855501626
